$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5 ("90min") raw data values B5:I5 ---
$ws.Range("B5").Value = 118.0
$ws.Range("C5").Value = 118.0
$ws.Range("D5").Value = 133.0
$ws.Range("E5").Value = 125.0
$ws.Range("F5").Value = 82.0
$ws.Range("G5").Value = 119.0
$ws.Range("H5").Value = 119.0
$ws.Range("I5").Value = 131.0

# Carry the formatting from row 2's raw-data cells down onto row 5's new cells
$ws.Range("B2:I2").Copy()
$ws.Range("B5:I5").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 5 summary formulas J5:M5 (extends the shared formulas from rows 2:4) ---
$ws.Range("J5").Formula = "=average(B5:E5)"
$ws.Range("K5").Formula = "=stdev(B5:E5)/sqrt(4)"
$ws.Range("L5").Formula = "=average(F5:I5)"
$ws.Range("M5").Formula = "=stdev(F5:I5)/sqrt(4)"

# Carry the formatting from row 4's formula cells down onto row 5's new formula cells
$ws.Range("J4:M4").Copy()
$ws.Range("J5:M5").PasteSpecial(-4122)  # xlPasteFormats

# --- New blank row 6, with A6 carrying the same row style as A5 ---
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A6").ClearContents()
